# Update gh-pages to output generated at 456a3b4
# Bumps several "想去人数" (F column) counts that were scraped from
# bilibili event listing pages across multiple sheets.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value = 2603
$ws1.Range("F7").Value = 97
$ws1.Range("F8").Value = 103
$ws1.Range("F20").Value = 7669
$ws1.Range("F21").Value = 8823
$ws1.Range("F32").Value = 1551

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 2675

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F9").Value = 2603
$ws4.Range("F11").Value = 97
$ws4.Range("F12").Value = 103
$ws4.Range("F24").Value = 7669
$ws4.Range("F25").Value = 7669
$ws4.Range("F26").Value = 8823
$ws4.Range("F34").Value = 1551
